# Adiamento do trabalho de lab. para sabado
# Moves the delivery date from 05/07 to 08/07:
#  1) In the "Data de entrega final:" line, strike through "05/07" and
#     append a bold, red "08/07" after it.
#  2) In the "(05/07)" parenthetical later in the text, turn it into
#     "(08/07)" by changing the "5" to "8" and re-marking the (moved)
#     hidden "_GoBack" bookmark right after the edited digit - this also
#     removes it from its old location at the very end of the document,
#     since Word only ever keeps a single "_GoBack" bookmark.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Hunk 1: "Data de entrega final: 05/07" -> strike "05/07" + " 08/07"
# ---------------------------------------------------------------------

$rng = $d.Content
$rng.Start = 0
$rng.Find.Execute("05/07", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$dateStart = $rng.Start
$dateEnd = $rng.End

# Duplicate the "05/07" run (keeps its sz/szCs) right after itself, then
# strike the original and turn the duplicate into " 08/07".
$rng.Copy()
$d.Range($dateEnd, $dateEnd).Paste()

$d.Range($dateStart, $dateEnd).Font.StrikeThrough = $true

$dupRange = $d.Range($dateEnd, $dateEnd + 5)
$dupRange.Text = " 08/07"

$newDateRange = $d.Range($dateEnd + 1, $dateEnd + 6)
$newDateRange.Font.Bold = $true
$newDateRange.Font.Color = 255

# ---------------------------------------------------------------------
# Hunk 2: "... fase (05/07), o grupo ..." -> "... fase (08/07), o grupo ..."
# split into "0" / "8" / bookmark "_GoBack" / "/07"
# ---------------------------------------------------------------------

$rng2 = $d.Content
$rng2.Start = $dateEnd + 6
$rng2.Find.Execute("05/07", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$secondStart = $rng2.Start

# "5" -> "8"
$d.Range($secondStart + 1, $secondStart + 2).Text = "8"

# Re-split "0" away from "8/07" without letting same-formatted runs
# re-merge: copy it, delete it, then paste it back in place.
$zeroRange = $d.Range($secondStart, $secondStart + 1)
$zeroRange.Copy()
$zeroRange.Delete()
$d.Range($secondStart, $secondStart).Paste()

# Re-seat (move) the hidden "_GoBack" bookmark between "8" and "/07".
# Word keeps only one "_GoBack" bookmark, so this also removes it from
# its previous location at the end of the document.
$bmRange = $d.Range($secondStart + 2, $secondStart + 2)
$d.Bookmarks.Add("_GoBack", $bmRange)
